$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newLabel = "Diferença 2025/07 - 2024/07"

# Row 2: Acre
$ws.Range("A2").Value = "Acre"
$ws.Range("B2").Value = $newLabel
$ws.Range("C2").Value = 2.06
$ws.Range("D2").Value = "1º"

# Row 3: Bahia
$ws.Range("A3").Value = "Bahia"
$ws.Range("B3").Value = $newLabel
$ws.Range("C3").Value = 1.78
$ws.Range("D3").Value = "2º"

# Row 4: Amapá
$ws.Range("A4").Value = "Amapá"
$ws.Range("B4").Value = $newLabel
$ws.Range("C4").Value = 1.77
$ws.Range("D4").Value = "3º"

# Row 5: Tocantins
$ws.Range("A5").Value = "Tocantins"
$ws.Range("B5").Value = $newLabel
$ws.Range("C5").Value = 1.74
$ws.Range("D5").Value = "4º"

# Row 6: Amazonas
$ws.Range("A6").Value = "Amazonas"
$ws.Range("B6").Value = $newLabel
$ws.Range("C6").Value = 1.5
$ws.Range("D6").Value = "5º"

# Row 7: Paraná
$ws.Range("A7").Value = "Paraná"
$ws.Range("B7").Value = $newLabel
$ws.Range("C7").Value = 1.07
$ws.Range("D7").Value = "6º"

# Row 8: Sergipe (name unchanged)
$ws.Range("B8").Value = $newLabel
$ws.Range("C8").Value = -3.77
$ws.Range("D8").Value = "27º"

# Row 9: Brasil (name unchanged)
$ws.Range("B9").Value = $newLabel
$ws.Range("C9").Value = 0.07000000000000001

# Row 10: Nordeste (name unchanged)
$ws.Range("B10").Value = $newLabel
$ws.Range("C10").Value = 0.29
